# Project_task_planning.xlsx update
# - Shorten W6 remark, split old W6/"Report" placeholder row into a real
#   W7 row (finish model evaluation / write report), and turn the old
#   W8/"Presentation" row's date into a real date value while recoloring
#   the highlighted (orange) row from the theme color to an explicit
#   orange (FFC000).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Row 7 (W6): shorten the note text, formatting untouched.
$ws.Range("C7").Value = "Both: finish feature selection"

# 2) Row 8 used to be a duplicate "W6" row whose only purpose was to hold
#    the word "Report" (yellow-highlighted). Turn it into the real W7 row,
#    keeping its existing (yellow) formatting untouched.
$ws.Range("A8").Value = "W7"
$ws.Range("B8").Value = "31/10/2022 6/11/2022"
$ws.Range("C8").Value = "Finish model evaluation, write report"

# 3) Row 9 (W8 / Presentation): turn the date range text into an actual
#    date value (11 Nov 2022 -> serial 44876), keep its "d-mmm" number
#    format which was already applied to that cell.
$ws.Range("B9").Value = 44876

# 4) Recolor the highlighted ("orange") fill used by row 9 from the old
#    theme color to the new explicit orange RGB(255,192,0) = FFC000.
#    (COM color longs are packed as 0x00BBGGRR -> 0 + (192*256) + 255 = 49407)
$ws.Range("A9:C9").Interior.Color = 49407

# 5) Update the active selection like in the saved workbook (C15 instead
#    of C16).
$ws.Range("C15").Select()
